# Update the hand-landmark (x, y) coordinates on the sheet for rows 2-79
# (landmark indices 0-77) as produced by the re-run / merge of
# hand_landmark_estimation.py. Column A holds the landmark index, column B
# the x coordinate and column C the y coordinate. The data set now also
# contains 14 additional landmarks (indices 64-77) that were not present
# before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 78,3
$arr[0,0] = 0
$arr[0,1] = 321.9075393676758
$arr[0,2] = 156.5245199203491
$arr[1,0] = 1
$arr[1,1] = 321.1606979370117
$arr[1,2] = 156.6767406463623
$arr[2,0] = 2
$arr[2,1] = 320.3333282470703
$arr[2,2] = 156.4820194244385
$arr[3,0] = 3
$arr[3,1] = 320.5591583251953
$arr[3,2] = 157.2189903259277
$arr[4,0] = 4
$arr[4,1] = 320.1755523681641
$arr[4,2] = 157.4060726165771
$arr[5,0] = 5
$arr[5,1] = 320.0727081298828
$arr[5,2] = 157.2195911407471
$arr[6,0] = 6
$arr[6,1] = 320.9890747070312
$arr[6,2] = 157.6928758621216
$arr[7,0] = 7
$arr[7,1] = 322.1222686767578
$arr[7,2] = 158.0412340164185
$arr[8,0] = 8
$arr[8,1] = 324.2362976074219
$arr[8,2] = 158.3535003662109
$arr[9,0] = 9
$arr[9,1] = 324.840087890625
$arr[9,2] = 157.8266572952271
$arr[10,0] = 10
$arr[10,1] = 325.7466125488281
$arr[10,2] = 157.1883773803711
$arr[11,0] = 11
$arr[11,1] = 325.0298309326172
$arr[11,2] = 156.2177610397339
$arr[12,0] = 12
$arr[12,1] = 321.9512176513672
$arr[12,2] = 154.4462299346924
$arr[13,0] = 13
$arr[13,1] = 317.0373725891113
$arr[13,2] = 149.9756526947021
$arr[14,0] = 14
$arr[14,1] = 308.841438293457
$arr[14,2] = 144.9496078491211
$arr[15,0] = 15
$arr[15,1] = 296.7984390258789
$arr[15,2] = 138.9054107666016
$arr[16,0] = 16
$arr[16,1] = 283.856201171875
$arr[16,2] = 133.2048225402832
$arr[17,0] = 17
$arr[17,1] = 268.5635757446289
$arr[17,2] = 127.4671411514282
$arr[18,0] = 18
$arr[18,1] = 253.7520217895508
$arr[18,2] = 124.6163177490234
$arr[19,0] = 19
$arr[19,1] = 237.3663520812988
$arr[19,2] = 123.1806135177612
$arr[20,0] = 20
$arr[20,1] = 221.2969398498535
$arr[20,2] = 121.8217420578003
$arr[21,0] = 21
$arr[21,1] = 206.0052490234375
$arr[21,2] = 124.066801071167
$arr[22,0] = 22
$arr[22,1] = 194.1226577758789
$arr[22,2] = 125.2685308456421
$arr[23,0] = 23
$arr[23,1] = 184.9114990234375
$arr[23,2] = 130.3762722015381
$arr[24,0] = 24
$arr[24,1] = 178.7818336486816
$arr[24,2] = 140.0556564331055
$arr[25,0] = 25
$arr[25,1] = 176.1171340942383
$arr[25,2] = 148.4401702880859
$arr[26,0] = 26
$arr[26,1] = 176.1407661437988
$arr[26,2] = 157.6689577102661
$arr[27,0] = 27
$arr[27,1] = 180.8287239074707
$arr[27,2] = 167.3396301269531
$arr[28,0] = 28
$arr[28,1] = 187.3880386352539
$arr[28,2] = 177.5644254684448
$arr[29,0] = 29
$arr[29,1] = 198.9091491699219
$arr[29,2] = 187.0843505859375
$arr[30,0] = 30
$arr[30,1] = 211.5151977539062
$arr[30,2] = 197.3724746704102
$arr[31,0] = 31
$arr[31,1] = 228.3205032348633
$arr[31,2] = 208.0852317810059
$arr[32,0] = 32
$arr[32,1] = 245.4554557800293
$arr[32,2] = 219.0063285827637
$arr[33,0] = 33
$arr[33,1] = 264.9747276306152
$arr[33,2] = 229.1789531707764
$arr[34,0] = 34
$arr[34,1] = 282.9477882385254
$arr[34,2] = 236.1346292495728
$arr[35,0] = 35
$arr[35,1] = 298.5334205627441
$arr[35,2] = 242.3202037811279
$arr[36,0] = 36
$arr[36,1] = 298.5334205627441
$arr[36,2] = 242.3202037811279
$arr[37,0] = 37
$arr[37,1] = 312.2934913635254
$arr[37,2] = 245.6957530975342
$arr[38,0] = 38
$arr[38,1] = 324.3391799926758
$arr[38,2] = 248.0888843536377
$arr[39,0] = 39
$arr[39,1] = 326.2100219726562
$arr[39,2] = 248.166675567627
$arr[40,0] = 40
$arr[40,1] = 323.5784530639648
$arr[40,2] = 246.4389324188232
$arr[41,0] = 41
$arr[41,1] = 318.4364700317383
$arr[41,2] = 246.1782360076904
$arr[42,0] = 42
$arr[42,1] = 310.3375625610352
$arr[42,2] = 245.2166175842285
$arr[43,0] = 43
$arr[43,1] = 299.0495109558105
$arr[43,2] = 242.4559020996094
$arr[44,0] = 44
$arr[44,1] = 287.0227241516113
$arr[44,2] = 240.7041263580322
$arr[45,0] = 45
$arr[45,1] = 272.3739814758301
$arr[45,2] = 240.2332878112793
$arr[46,0] = 46
$arr[46,1] = 272.3739814758301
$arr[46,2] = 240.2332878112793
$arr[47,0] = 47
$arr[47,1] = 258.166618347168
$arr[47,2] = 241.8719673156738
$arr[48,0] = 48
$arr[48,1] = 243.9486694335938
$arr[48,2] = 246.6483306884766
$arr[49,0] = 49
$arr[49,1] = 229.3948554992676
$arr[49,2] = 249.8833465576172
$arr[50,0] = 50
$arr[50,1] = 209.7288131713867
$arr[50,2] = 261.0505485534668
$arr[51,0] = 51
$arr[51,1] = 202.0956611633301
$arr[51,2] = 269.9016380310059
$arr[52,0] = 52
$arr[52,1] = 196.9227027893066
$arr[52,2] = 277.631950378418
$arr[53,0] = 53
$arr[53,1] = 193.6286163330078
$arr[53,2] = 286.9873809814453
$arr[54,0] = 54
$arr[54,1] = 191.7404556274414
$arr[54,2] = 295.532283782959
$arr[55,0] = 55
$arr[55,1] = 191.7404556274414
$arr[55,2] = 295.532283782959
$arr[56,0] = 56
$arr[56,1] = 195.1507568359375
$arr[56,2] = 303.2415390014648
$arr[57,0] = 57
$arr[57,1] = 201.2717247009277
$arr[57,2] = 310.9980583190918
$arr[58,0] = 58
$arr[58,1] = 211.2743377685547
$arr[58,2] = 317.8551578521729
$arr[59,0] = 59
$arr[59,1] = 224.1765594482422
$arr[59,2] = 326.4860343933105
$arr[60,0] = 60
$arr[60,1] = 237.4317359924316
$arr[60,2] = 332.1601295471191
$arr[61,0] = 61
$arr[61,1] = 253.2627487182617
$arr[61,2] = 337.410192489624
$arr[62,0] = 62
$arr[62,1] = 288.6498641967773
$arr[62,2] = 345.6345748901367
$arr[63,0] = 63
$arr[63,1] = 304.8026657104492
$arr[63,2] = 344.3311214447021
$arr[64,0] = 64
$arr[64,1] = 304.8026657104492
$arr[64,2] = 344.3311214447021
$arr[65,0] = 65
$arr[65,1] = 322.8125381469727
$arr[65,2] = 346.3611316680908
$arr[66,0] = 66
$arr[66,1] = 336.7195129394531
$arr[66,2] = 345.2759456634521
$arr[67,0] = 67
$arr[67,1] = 351.0630416870117
$arr[67,2] = 347.7375984191895
$arr[68,0] = 68
$arr[68,1] = 363.4733963012695
$arr[68,2] = 346.0503959655762
$arr[69,0] = 69
$arr[69,1] = 371.406135559082
$arr[69,2] = 343.9895439147949
$arr[70,0] = 70
$arr[70,1] = 377.1853637695312
$arr[70,2] = 342.2107315063477
$arr[71,0] = 71
$arr[71,1] = 380.9870529174805
$arr[71,2] = 340.9621524810791
$arr[72,0] = 72
$arr[72,1] = 381.7733764648438
$arr[72,2] = 341.1250591278076
$arr[73,0] = 73
$arr[73,1] = 380.7618713378906
$arr[73,2] = 362.5720882415771
$arr[74,0] = 74
$arr[74,1] = 375.5272674560547
$arr[74,2] = 388.2195568084717
$arr[75,0] = 75
$arr[75,1] = 363.931770324707
$arr[75,2] = 400.9485912322998
$arr[76,0] = 76
$arr[76,1] = 360.24658203125
$arr[76,2] = 405.137300491333
$arr[77,0] = 77
$arr[77,1] = 358.3580017089844
$arr[77,2] = 407.5924301147461

$ws.Range("A2:C79").Value = $arr
